$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = "Riomaior"
$ws.Range("B16").Value = "Rio"

$ws.Range("C18").Select()
